$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "2021" column (R) by cloning the formatting of the existing
# "2020" column (Q) cell-by-cell, then filling in the new values.

# Row 2: blank separator cell (same style/border as Q2, no value)
$ws.Range("Q2").Copy($ws.Range("R2"))
$ws.Range("R2").Value = ""

# Row 3: year header
$ws.Range("Q3").Copy($ws.Range("R3"))
$ws.Range("R3").Value = 2021

# Row 4: population count
$ws.Range("Q4").Copy($ws.Range("R4"))
$ws.Range("R4").Value = 202551

# Row 5: percentage of total population
$ws.Range("Q5").Copy($ws.Range("R5"))
$ws.Range("R5").Value = 2.9794303052841493

# Move the active selection to the newly added column, as in the authored
# workbook.
[void]$ws.Range("R2").Select()
